# MHD2-259: Report template and related changes for reporting on 136 genes
#
# This script applies the textual edits captured in the commit diff:
#  1. Merge the runs (and drop the interleaved spell-check proofErr markers)
#     in the "DNA is analysed..." / "A custom pipeline..." / "Variants are
#     analysed..." method paragraph - text content is unchanged, only the
#     run/proofErr bookkeeping collapses, which Word does naturally when a
#     Find/Replace spans several runs.
#  2. Same run-merge cleanup for the "subcategorisation" sentence.
#  3. Same run-merge cleanup for the JAK2/ASXL1 detection-limit sentence.
#  4. Same proofErr cleanup for the lone "LIMITATIONS_cfDNA_IN" run.
#  5. Remove the trailing "; please note FLT3-ITDs and UBTF-TDs ... if
#     sample tested." sentence from the CDS footnote paragraph.
#  6. Bump the SAVEDATE field result text from 5-Nov-2025 to 7-Nov-2025.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "WARNING: text not found -> $old"
    }
    return $result
}

# --- 1. Method paragraph: merge runs / drop proofErr, in three chunks so the
#        existing bookmark (_Hlk177382982) around the Oncoanalyser sentence
#        stays exactly where it is. Replacement text == search text (the
#        visible wording does not change).

$seg1 = "DNA is analysed by targeted gene sequencing of coding regions and flanking splice sites (within 2 bp) of the genes listed below. Libraries are prepared using a custom Twist Bioscience target enrichment panel (Peter MacCallum Cancer Centre AllHaem DNA Twist v2, design ID TE-91041418) and sequenced on an Illumina NovaSeq X Plus (Australian Genome Research Facility) with 150 bp paired end reads. "
Replace-Text $seg1 $seg1

$seg2 = "A custom pipeline utilising the Oncoanalyser analysis pipeline (OncoPath v1) is used to generate aligned reads and call variants (single nucleotide variants and short insertions or deletions) against the hg19 human reference genome. "
Replace-Text $seg2 $seg2

$seg3 = "Variants are analysed using PathOS software (Peter Mac) and described according to HGVS nomenclature version 19.01 (http://varnomen.hgvs.org/) with minor differences in accordance with Peter MacCallum Cancer Centre Molecular Pathology departmental policy. The following population variation and cancer or genetic disease databases are commonly used in addition to literature review to assist with variant interpretation: the Genome Aggregation Database (gnomAD; gnomad.broadinstitute.org), the Catalogue of Somatic Mutations in Cancer (COSMIC; cancer.sanger.ac.uk), ClinVar (ncbi.nlm.nih.gov/clinvar) and the IARC TP53 Database (p53.iarc.fr). "
Replace-Text $seg3 $seg3

# --- 2. "subcategorisation" sentence: merge runs / drop proofErr.
$seg4 = " (the variant either defines a diagnostic category or is sufficiently specific for the clinical context to contribute to diagnostic subcategorisation), "
Replace-Text $seg4 $seg4

# --- 3. JAK2 / ASXL1 detection limit sentence: merge runs / drop proofErr.
$seg5 = "The detection limit of this assay for specimens sequenced to the target read depth of 500x is a variant allele frequency (VAF) of approximately 2% with the exception of JAK2 c.1849G>T;p.(Val617Phe) (detection limit ~ 1%) and ASXL1 c.1934dup;p.(Gly646Trpfs*12) (detection limit ~ 5%). This assay is primarily qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. Insertions or deletions (particularly those > 25 bp in length or in homopolymer regions), including FLT3-ITDs and UBTF-TDs, are not reliably detected by this assay. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. This assay does not distinguish between somatic and germline variants. In addition, the clonal origin of somatic variants (i.e. disease compartment or cell lineage) cannot be determined. "
Replace-Text $seg5 $seg5

# --- 4. Lone "LIMITATIONS_cfDNA_IN" heading run: drop surrounding proofErr.
$seg6 = "LIMITATIONS_cfDNA_IN"
Replace-Text $seg6 $seg6

# --- 5. Trim the CDS footnote paragraph back to just the abbreviation key.
$seg7 = "; please note FLT3-ITDs and UBTF-TDs are not reliably detected with this assay. A separate assay may have been performed, result included in Test Description if sample tested."
Replace-Text $seg7 ""

# --- 6. Bump the report SAVEDATE field result.
Replace-Text "5-Nov-2025" "7-Nov-2025"
